$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 2.6
$ws.Range("I3").Value = 3.1
$ws.Range("K3").Value = 1.91
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6.5
$ws.Range("AJ3").Value = 34
$ws.Range("AN3").Value = 4.33
$ws.Range("AX3").Value = 19
$ws.Range("AZ3").Value = 67
